$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0.274, 0.601, 0.471, -82, 3, 65),
    @(0.323, 0.361, 0.955, -126, 3, -8),
    @(-0.076, -0.317, 0.499, -42, 28, 173),
    @(0.029, 0.194, 0.108, 87, -28, -84),
    @(0.04, -0.184, 0.331, -77, 37, -129),
    @(0.25, -0.461, 0.286, 76, 21, -87),
    @(0.07199999999999999, -0.228, 0.027, -15, -45, 80),
    @(-0.149, 0.178, 0.671, 51, 0, 172),
    @(0.747, -0.077, 0.675, 53, 25, -160),
    @(-0.092, 0.048, 0.436, -51, -54, 100)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $col = $j + 1
        $ws.Cells.Item($row, $col).Value = $rowValues[$j]
    }
}
